$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8172028064727783
$ws.Range("B1").Value = 1.182101249694824
$ws.Range("C1").Value = 3.56993842124939
$ws.Range("D1").Value = 3.956349611282349
$ws.Range("E1").Value = 1.233461260795593
